# Matrix of TCAT.xlsx - "Add files via upload" edit
# Extends Sheet1's matrix from 10x10 to 15x15, tweaks one existing cell
# (E2: 3 -> 4), and adds a brand-new "Sheet2" with a 5x5 matrix.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: full 15x15 matrix (rows 1-15, cols A-O) -----------------------
$sheet1Rows = @(
    @(3,0,0,1,0,0,0,0,1,0,0,0,0,0,0),
    @(0,5,3,5,4,0,0,2,6,1,0,0,0,0,1),
    @(0,3,5,1,1,0,0,0,1,0,0,0,0,1,0),
    @(1,5,1,13,3,1,2,3,5,1,2,3,1,1,3),
    @(0,4,1,3,5,0,0,0,2,0,0,0,0,1,1),
    @(0,0,0,1,0,2,0,0,0,0,0,1,0,0,0),
    @(0,0,0,2,0,0,1,0,0,0,0,0,0,0,0),
    @(0,2,0,3,0,0,0,6,2,0,0,1,0,1,3),
    @(1,6,1,5,2,0,0,2,7,0,0,0,0,2,2),
    @(0,1,0,1,0,0,0,0,0,1,0,0,0,0,0),
    @(0,0,0,2,0,0,0,0,0,0,2,0,2,0,0),
    @(0,0,0,3,0,1,0,1,0,0,0,3,0,0,1),
    @(0,0,0,1,0,0,0,0,0,0,2,0,2,0,0),
    @(0,0,1,1,1,0,0,1,2,0,0,0,0,4,1),
    @(0,1,0,3,1,0,0,3,2,0,0,1,0,1,1)
)

for ($r = 0; $r -lt $sheet1Rows.Length; $r++) {
    $rowData = $sheet1Rows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws1.Cells.Item($r + 1, $c + 1).Value = $rowData[$c]
    }
}

# Selection left on Sheet1 after the edit
[void]$ws1.Range("D4").Select()

# --- Add Sheet2 right after Sheet1 -----------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$sheet2Rows = @(
    @(2,0,2,0,0),
    @(0,3,0,0,0),
    @(2,0,2,0,0),
    @(0,0,0,4,0),
    @(0,0,0,0,1)
)

for ($r = 0; $r -lt $sheet2Rows.Length; $r++) {
    $rowData = $sheet2Rows[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws2.Cells.Item($r + 1, $c + 1).Value = $rowData[$c]
    }
}

[void]$ws2.Range("D3").Select()

# Leave Sheet1 as the active sheet/tab (tabSelected="1" in target)
$ws1.Activate()
